$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Formula = "'" + '62.874.22'
$ws.Range('E2').Formula = '  +4.83%  '
$ws.Range('D3').Formula = "'" + '3.357.26'
$ws.Range('E3').Formula = '  +5.11%  '
$ws.Range('E4').Formula = '  -0.01%  '
$ws.Range('D5').Formula = "'" + '556.82'
$ws.Range('E5').Formula = '  +3.47%  '
$ws.Range('D6').Formula = "'" + '153.63'
$ws.Range('E6').Formula = '  +6.04%  '
$ws.Range('E7').Formula = '  -0.10%  '
$ws.Range('E8').Formula = '  +0.64%  '
$ws.Range('E9').Formula = '  +2.48%  '
$ws.Range('E10').Formula = '  +4.48%  '
$ws.Range('E11').Formula = '  +1.91%  '
$ws.Range('D12').Formula = "'" + '3.933.89'
$ws.Range('E12').Formula = '  +5.00%  '
$ws.Range('D13').Formula = "'" + '0.138'
$ws.Range('E13').Formula = '  +0.41%  '
$ws.Range('B14').Formula = 'ShibaInu'
$ws.Range('C14').Formula = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D14').Formula = "'" + '0.0000182'
$ws.Range('E14').Formula = '  +3.95%  '
$ws.Range('B15').Formula = 'Avalanche'
$ws.Range('C15').Formula = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D15').Formula = "'" + '27.11'
$ws.Range('E15').Formula = '  +4.27%  '
$ws.Range('D16').Formula = "'" + '62.942.75'
$ws.Range('E16').Formula = '  +4.86%  '
$ws.Range('D17').Formula = "'" + '3.352.49'
$ws.Range('E17').Formula = '  +4.92%  '
$ws.Range('D18').Formula = "'" + '6.51'
$ws.Range('E18').Formula = '  +4.64%  '
$ws.Range('E19').Formula = '  +5.49%  '
$ws.Range('D20').Formula = "'" + '8.45'
$ws.Range('E20').Formula = '  +1.20%  '
$ws.Range('D21').Formula = "'" + '388.42'
$ws.Range('E21').Formula = '  +1.42%  '
$ws.Range('E22').Formula = '  +2.40%  '
$ws.Range('D23').Formula = "'" + '1.00'
$ws.Range('E23').Formula = '  +0.15%  '
$ws.Range('D24').Formula = "'" + '70.52'
$ws.Range('E24').Formula = '  +0.33%  '
$ws.Range('D25').Formula = "'" + '0.180'
$ws.Range('E25').Formula = '  +4.93%  '
$ws.Range('D26').Formula = "'" + '8.90'
$ws.Range('E26').Formula = '  +0.31%  '
$ws.Range('D27').Formula = "'" + '0.0₃0975'
$ws.Range('E27').Formula = '  +7.93%  '
$ws.Range('E28').Formula = '  +0.40%  '
$ws.Range('D29').Formula = "'" + '6.65'
$ws.Range('E29').Formula = '  +7.42%  '
$ws.Range('E30').Formula = '  +4.36%  '
$ws.Range('D31').Formula = "'" + '5.66'
$ws.Range('E31').Formula = '  +5.42%  '
$ws.Range('D32').Formula = "'" + '23.08'
$ws.Range('E32').Formula = '  +2.94%  '
$ws.Range('E33').Formula = '  +7.12%  '
$ws.Range('E34').Formula = '  +0.86%  '
$ws.Range('D35').Formula = "'" + '160.53'
$ws.Range('E35').Formula = '  +2.70%  '
$ws.Range('E36').Formula = '  +8.97%  '
$ws.Range('D37').Formula = "'" + '1.90'
$ws.Range('E37').Formula = '  +12.44%  '
$ws.Range('D38').Formula = "'" + '27.01'
$ws.Range('E38').Formula = '  +4.78%  '
$ws.Range('D39').Formula = "'" + '0.0742'
$ws.Range('E39').Formula = '  +4.19%  '
$ws.Range('D40').Formula = "'" + '2.831.96'
$ws.Range('E40').Formula = '  +1.73%  '
$ws.Range('D41').Formula = "'" + '0.0311'
$ws.Range('E41').Formula = '  +8.58%  '
$ws.Range('D42').Formula = "'" + '4.33'
$ws.Range('E42').Formula = '  +1.87%  '
$ws.Range('B43').Formula = 'Mantle'
$ws.Range('C43').Formula = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D43').Formula = "'" + '0.749'
$ws.Range('E43').Formula = '  +2.70%  '
$ws.Range('B44').Formula = 'OKB'
$ws.Range('C44').Formula = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D44').Formula = "'" + '40.75'
$ws.Range('E44').Formula = '  +2.53%  '
$ws.Range('E45').Formula = '  +5.05%  '
$ws.Range('D46').Formula = "'" + '22.19'
$ws.Range('E46').Formula = '  +7.92%  '
$ws.Range('D47').Formula = "'" + '3.399.77'
$ws.Range('E47').Formula = '  +5.07%  '
$ws.Range('D48').Formula = "'" + '0.103'
$ws.Range('E48').Formula = '  +2.19%  '
$ws.Range('D49').Formula = "'" + '6.32'
$ws.Range('E49').Formula = '  +2.18%  '
$ws.Range('D50').Formula = "'" + '0.809'
$ws.Range('E50').Formula = '  +0.89%  '
$ws.Range('D51').Formula = "'" + '280.73'
$ws.Range('E51').Formula = '  +6.17%  '
